# Apply the edit described by the diff:
# 1. Update the cached "last updated" date field text from 3/23/2015 to
#    3/24/2015 everywhere it is cached (slide master, all slide layouts,
#    and the notes master).
# 2. On slide 1, delete the "TextBox 7" shape (the motivational-quote
#    placeholder textbox) along with its entrance animation effect and
#    its corresponding build (bldLst) entry.

$p = $ppt.ActivePresentation

function Update-DateFields($shapes) {
    foreach ($shp in $shapes) {
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "3/23/2015") {
                $tr.Text = "3/24/2015"
            }
        }
    }
}

# --- slide master ---
Update-DateFields $p.SlideMaster.Shapes

# --- every slide layout ---
foreach ($layout in $p.SlideMaster.CustomLayouts) {
    Update-DateFields $layout.Shapes
}

# --- notes master ---
Update-DateFields $p.NotesMaster.Shapes

# --- slide 1: remove the "TextBox 7" placeholder textbox & its animation ---
$s1 = $p.Slides.Item(1)

# Remove the entrance animation effect targeting the textbox, if present.
$seq = $s1.TimeLine.MainSequence
for ($i = $seq.Count; $i -ge 1; $i--) {
    $eff = $seq.Item($i)
    if ($eff.Shape.Name -eq "TextBox 7") {
        $eff.Delete()
    }
}

# Remove the shape itself (this also clears its build entry).
foreach ($shp in $s1.Shapes) {
    if ($shp.Name -eq "TextBox 7") {
        $shp.Delete()
    }
}
